# Refresh the NATMI ligand-receptor table (Dkk2 -> Kremen1) with numbers
# recomputed from TPM instead of raw counts.
#
# The old sheet only had a "FAPs" sending cluster (3 rows, one per target
# cluster: ECs / FAPs / MuSCs). The TPM re-run adds a second sending
# cluster, "ECs", with its own 3 rows, while recomputing every numeric
# column (columns E:T) for all six sending/target combinations. The
# "FAPs" sending-cluster rows keep their place in the sheet but move down
# to rows 5-7 and pick up refreshed specificity figures (S:T) as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol | D Target cluster
#          E..T the 16 numeric NATMI metrics
$data = @(
    @("ECs",  "Dkk2", "Kremen1", "ECs",
      1, 0.3333333333333333, 0.1293116666666667, 0.387935,
      0.0151634279701504, 0.0151634279701504, 3, 1,
      3.636278, 10.908834, 0.1254799416243108, 0.1254799416243107,
      0.4702131686433333, 4.23191851779, 0.001902706056518913, 0.001902706056518913),
    @("ECs",  "Dkk2", "Kremen1", "FAPs",
      1, 0.3333333333333333, 0.1293116666666667, 0.387935,
      0.0151634279701504, 0.0151634279701504, 3, 1,
      17.985402, 53.956206, 0.620636594080475, 0.620636594080475,
      2.32572230829, 20.93150077461, 0.009410978289978755, 0.009410978289978757),
    @("ECs",  "Dkk2", "Kremen1", "MuSCs",
      1, 0.3333333333333333, 0.1293116666666667, 0.387935,
      0.0151634279701504, 0.0151634279701504, 3, 1,
      7.357278333333333, 22.071835, 0.2538834642952142, 0.2538834642952141,
      0.9513819234138888, 8.562437310724999, 0.003849743623652732, 0.003849743623652731),
    @("FAPs", "Dkk2", "Kremen1", "ECs",
      3, 1, 8.398553333333334, 25.19566,
      0.9848365720298496, 0.9848365720298496, 3, 1,
      3.636278, 10.908834, 0.1254799416243108, 0.1254799416243107,
      30.53947471782667, 274.85527246044, 0.1235772355677918, 0.1235772355677918),
    @("FAPs", "Dkk2", "Kremen1", "FAPs",
      3, 1, 8.398553333333334, 25.19566,
      0.9848365720298496, 0.9848365720298496, 3, 1,
      17.985402, 53.956206, 0.620636594080475, 0.620636594080475,
      151.05135791844, 1359.46222126596, 0.6112256157904963, 0.6112256157904963),
    @("FAPs", "Dkk2", "Kremen1", "MuSCs",
      3, 1, 8.398553333333334, 25.19566,
      0.9848365720298496, 0.9848365720298496, 3, 1,
      7.357278333333333, 22.071835, 0.2538834642952142, 0.2538834642952141,
      61.79049447067779, 556.1144502361, 0.2500337206715615, 0.2500337206715614)
)

$row = 2
foreach ($values in $data) {
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
    $row++
}
